$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = "25.657.60"
$ws.Range("E2").Value = "  -5.79%  "

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = "1.809.71"
$ws.Range("E3").Value = "  -4.99%  "

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").Value = "'276.65"
$ws.Range("E5").Value = "  -9.71%  "

# Row 6: 'USDC' -> 'USDC'
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.00%  "

# Row 7: 'XRP' -> 'XRP'
$ws.Range("D7").Value = "'0.5015"
$ws.Range("E7").Value = "  -6.31%  "

# Row 8: 'Cardano' -> 'Cardano'
$ws.Range("D8").Value = "'0.3500"
$ws.Range("E8").Value = "  -8.38%  "

# Row 9: 'OKB' -> 'Dogecoin'
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.06625"
$ws.Range("E9").Value = "  -9.26%  "

# Row 10: 'Dogecoin' -> 'Solana'
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "'20.05"
$ws.Range("E10").Value = "  -9.92%  "

# Row 11: 'Solana' -> 'Polygon'
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'0.8381"
$ws.Range("E11").Value = "  -7.50%  "

# Row 12: 'Polygon' -> 'TRON'
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07822"
$ws.Range("E12").Value = "  -4.80%  "

# Row 13: 'TRON' -> 'WrappedEther'
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.811.80"
$ws.Range("E13").Value = "  +71.30%  "

# Row 14: 'WrappedEther' -> 'Polkadot'
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.034"
$ws.Range("E14").Value = "  -5.91%  "

# Row 15: 'Polkadot' -> 'Litecoin'
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'87.35"
$ws.Range("E15").Value = "  -8.75%  "

# Row 16: 'Litecoin' -> 'BinanceUSD'
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.11%  "

# Row 17: 'BinanceUSD' -> 'Avalanche'
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'13.87"
$ws.Range("E17").Value = "  -6.67%  "

# Row 18: 'Avalanche' -> 'Dai'
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.00%  "

# Row 19: 'Dai' -> 'ShibaInu'
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007942"
$ws.Range("E19").Value = "  -8.41%  "

# Row 20: 'ShibaInu' -> 'WrappedBTC'
$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "25.726.37"
$ws.Range("E20").Value = "  -5.62%  "

# Row 21: 'WrappedBTC' -> 'Uniswap'
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'4.709"
$ws.Range("E21").Value = "  -6.80%  "

# Row 22: 'Uniswap' -> 'Cosmos'
$ws.Range("B22").Value = "Cosmos"
$ws.Range("C22").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D22").Value = "'9.984"
$ws.Range("E22").Value = "  -7.43%  "

# Row 23: 'Cosmos' -> 'Chainlink'
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "'6.058"
$ws.Range("E23").Value = "  -7.14%  "

# Row 24: 'Chainlink' -> 'Monero'
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").Value = "'141.22"
$ws.Range("E24").Value = "  -5.27%  "

# Row 25: 'Monero' -> 'Toncoin'
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'1.662"
$ws.Range("E25").Value = "  -4.75%  "

# Row 26: 'Toncoin' -> 'LidoDAOToken'
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.088"
$ws.Range("E26").Value = "  -9.02%  "

# Row 27: 'LidoDAOToken' -> 'EthereumClassic'
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'16.86"
$ws.Range("E27").Value = "  -8.30%  "

# Row 28: 'EthereumClassic' -> 'BitcoinCash'
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").Value = "'108.54"
$ws.Range("E28").Value = "  -7.13%  "

# Row 29: 'BitcoinCash' -> 'InternetComputer(DFINITY)'
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'4.296"
$ws.Range("E29").Value = "  -10.88%  "

# Row 30: 'InternetComputer(DFINITY)' -> 'Filecoin'
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'4.197"
$ws.Range("E30").Value = "  -11.23%  "

# Row 31: 'Filecoin' -> 'Stellar'
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.08798"
$ws.Range("E31").Value = "  -4.58%  "

# Row 32: 'Stellar' -> 'Hedera'
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.04825"
$ws.Range("E32").Value = "  -5.11%  "

# Row 33: 'Hedera' -> 'ImmutableX'
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.7379"
$ws.Range("E33").Value = "  -11.12%  "

# Row 34: 'ImmutableX' -> 'HuobiToken'
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.867"
$ws.Range("E34").Value = "  -4.62%  "

# Row 35: 'HuobiToken' -> 'ARBITRUM'
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.122"
$ws.Range("E35").Value = "  -7.92%  "

# Row 36: 'ARBITRUM' -> 'Frax'
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "'1.001"
$ws.Range("E36").Value = "  +0.07%  "

# Row 37: 'Frax' -> 'MXToken'
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'3.038"
$ws.Range("E37").Value = "  -8.43%  "

# Row 38: 'MXToken' -> 'RenderToken'
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.390"
$ws.Range("E38").Value = "  -10.84%  "

# Row 39: 'RenderToken' -> 'VeChain'
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01864"
$ws.Range("E39").Value = "  -6.96%  "

# Row 40: 'TheSandbox' -> 'TheSandbox'
$ws.Range("D40").Value = "'0.5190"
$ws.Range("E40").Value = "  -11.77%  "

# Row 41: 'VeChain' -> 'TrustWalletToken'
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.9667"
$ws.Range("E41").Value = "  -10.32%  "

# Row 42: 'TrustWalletToken' -> 'FraxShare'
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.201"
$ws.Range("E42").Value = "  -6.60%  "

# Row 43: 'FraxShare' -> 'Quant'
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'111.32"
$ws.Range("E43").Value = "  -5.13%  "

# Row 44: 'Quant' -> 'Aptos'
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'8.123"
$ws.Range("E44").Value = "  -13.18%  "

# Row 45: 'Aptos' -> 'PaxDollar'
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46: 'PaxDollar' -> 'Decentraland'
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.4608"
$ws.Range("E46").Value = "  -9.86%  "

# Row 47: 'Decentraland' -> 'Algorand'
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1379"
$ws.Range("E47").Value = "  -9.53%  "

# Row 48: 'Algorand' -> 'EnergySwap'
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.216"
$ws.Range("E48").Value = "  -9.22%  "

# Row 49: 'EnergySwap' -> 'Elrond'
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'35.64"
$ws.Range("E49").Value = "  -7.11%  "

# Row 50: 'Elrond' -> 'NEARProtocol'
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.487"
$ws.Range("E50").Value = "  -9.51%  "

# Row 51: 'NEARProtocol' -> 'Cronos'
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05852"
$ws.Range("E51").Value = "  -4.98%  "
